# Recovered Registration System from data corruption
# Re-adds the "last name" column (P1:P5) and restores a stray duplicate
# value (Q6), plus fixes a corrupted row-4 id value in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "last name" column entries (creates shared strings last1..last5)
$ws.Range("P1").Value = "last1"
$ws.Range("P2").Value = "last2"
$ws.Range("P3").Value = "last3"
$ws.Range("P4").Value = "last4"
$ws.Range("P5").Value = "last5"

# Row 4's id was corrupted to "4"; restore the real value
$ws.Range("A4").Value = 12524734

# Restore stray recovered value in Q6 (same text as B2/C2: "asd")
$ws.Range("Q6").Value = "asd"

# Leave the selection where the recovery work finished
$ws.Range("Q6").Select()
